# Add support for shared strings
#
# - A2  -> "String"   (new shared string, index 0)
# - B3  -> "Foo bar"  (new shared string, index 1)
# - F4  -> "Foo bar"  (re-uses shared string index 1)
# - Selection moves from D1 to D2
# - The sheet's used range grows from A1:D3 to A1:F4
# - The explicit per-column width metadata (<cols>) present in the source
#   file is dropped, which ClearFormats also achieves since none of the
#   cells in this sheet carry any real formatting of their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the workbook's inherited custom column-width metadata (the sheet has
# no real formatting, so this is a no-op for appearance but clears <cols>).
$ws.Cells.ClearFormats() | Out-Null

# New string-valued cells (these populate xl/sharedStrings.xml on save).
$ws.Range("A2").Value = "String"
$ws.Range("B3").Value = "Foo bar"
$ws.Range("F4").Value = "Foo bar"

# Match the new active selection.
$ws.Range("D2").Select()

# Cosmetic tab-ratio tweak recorded in the workbook view.
$excel.ActiveWindow.TabRatio = 0.992
